$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.148.76'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.882.83'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'313.59"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = "'0.5068"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('D8').Value = "'0.3857"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('D9').Value = "'0.09060"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').Value = "'6.373"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').Value = "'20.84"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '1.873.62'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').Value = "'7.284"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = "'0.00001114"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = "'91.49"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = "'0.06597"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = "'18.27"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.63%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').Value = "'6.133"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.22%  '
$ws.Range('D23').Value = '28.169.37'
$ws.Range('E23').Value = '  +0.40%  '
$ws.Range('D24').Value = "'11.48"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('D25').Value = "'2.265"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').Value = "'2.557"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('D27').Value = '2.089.74'
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('D28').Value = "'20.88"
$ws.Range('D28').ClearFormats()
$ws.Range('D29').Value = "'156.73"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('D30').Value = "'127.19"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').Value = "'1.065"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('D33').Value = "'5.632"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('D34').Value = "'3.597"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').Value = "'9.614"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').Value = "'0.06626"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').Value = "'0.02410"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').Value = "'0.2198"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').Value = "'1.292"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('D40').Value = "'1.216"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('E41').Value = '  +1.62%  '
$ws.Range('D42').Value = "'11.53"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').Value = "'4.937"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.02%  '
$ws.Range('D44').Value = "'0.6056"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.52%  '
$ws.Range('D45').Value = "'13.22"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = "'3.668"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').Value = "'1.275"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = "'1.243"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.75%  '
$ws.Range('D49').Value = "'2.009"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').Value = "'121.41"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('D51').Value = "'79.78"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.85%  '
